$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Repository URI"
$ws.Range("C8").Value = "Created at"
$ws.Range("C12").Value = "Counts"
